$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.828.84'
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = '1.708.30'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.11'
$ws.Range("E5").Value = '  +1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3754'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.65'
$ws.Range("E8").Value = '  +3.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3445'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.207'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07543'
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.10'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.296'
$ws.Range("E14").Value = '  +2.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.039'
$ws.Range("E15").Value = '  +4.06%  '
$ws.Range("D16").Value = '1.707.69'
$ws.Range("E16").Value = '  +1.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001135'
$ws.Range("E17").Value = '  +2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06726'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '84.42'
$ws.Range("E20").Value = '  +3.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.31'
$ws.Range("E21").Value = '  +5.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.376'
$ws.Range("E22").Value = '  +4.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.23'
$ws.Range("E23").Value = '  +10.41%  '
$ws.Range("D24").Value = '24.812.75'
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.455'
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.792'
$ws.Range("E26").Value = '  +3.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.38'
$ws.Range("E27").Value = '  +3.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.20'
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '132.54'
$ws.Range("E29").Value = '  +4.18%  '
$ws.Range("D30").Value = '1.898.38'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.242'
$ws.Range("E31").Value = '  +27.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.939'
$ws.Range("E32").Value = '  +9.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.223'
$ws.Range("E33").Value = '  +5.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.831'
$ws.Range("E34").Value = '  +5.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.77'
$ws.Range("E35").Value = '  +11.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08804'
$ws.Range("E36").Value = '  +3.49%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.334'
$ws.Range("E37").Value = '  +2.87%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.604'
$ws.Range("E38").Value = '  +4.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06674'
$ws.Range("E39").Value = '  +2.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02408'
$ws.Range("E40").Value = '  +2.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2236'
$ws.Range("E41").Value = '  +5.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.279'
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6446'
$ws.Range("E43").Value = '  +3.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.95'
$ws.Range("E45").Value = '  +6.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6165'
$ws.Range("E46").Value = '  +3.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.823'
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.136'
$ws.Range("E48").Value = '  +5.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.19'
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07312'
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.74'
$ws.Range("E51").Value = '  +5.08%  '
